# DataPath.xlsx update: "Added RCWA and TMM, fixed FDTD"
# - Fix an existing FDTD row (F11: 1 -> 0) on "Raw FDTD Data"
# - Add a new "Info file stub" column (U on sheet1, T on sheet2) with values
#   for all existing rows
# - Append several new simulation rows to "Raw FDTD Data" (rows 12-19) and
#   "Free Space Raw FDTD Data" (rows 4-15)
# - Add the Ex/Ey/Ez/Hx/Hy/Hz coefficient columns (N:S) to
#   "Free Space Raw FDTD Data" for existing rows 2-3
# - Restore the selections on both sheets

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Raw FDTD Data")
$ws2 = $wb.Worksheets.Item("Free Space Raw FDTD Data")

# ---------------------------------------------------------------------------
# Sheet "Raw FDTD Data"
# ---------------------------------------------------------------------------

# New header for column U
$ws1.Range("U1").Value = 'Info file stub'

# Existing rows 2-11 get the new "Info file stub" value
$ws1.Range("U2").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U3").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U4").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U5").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U6").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U7").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U8").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U9").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U10").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'
$ws1.Range("U11").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# Row 11's "dipole only" flag flips back off now that FDTD is fixed
$ws1.Range("F11").Value = 0

# New row 12
$ws1.Range("A12").Value = 'K:\Antropy\ARC waveguide 2_2 with slanted mirror one end long waveguide material below.xf\Simulations\000007\Run0001\output'
$ws1.Range("B12").Value = 'S2F design +x direction'
$ws1.Range("C12").Value = 2.2
$ws1.Range("D12").Value = 12000
$ws1.Range("E12").Value = [double]"1.7314999999999999E-16"
$ws1.Range("F12").Value = 0
$ws1.Range("G12").Value = 1000
$ws1.Range("H12").Value = 1050
$ws1.Range("I12").Value = 166
$ws1.Range("J12").Value = 2
$ws1.Range("K12").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response\Free Space 700 to 1200 nm with 5 nm mesh.mat'
$ws1.Range("L12").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response'
$ws1.Range("M12").Value = 'MultiPoint__X_sensor_4_transient_'
$ws1.Range("N12").Value = 'MultiPoint__X_sensor_4_transient_'
$ws1.Range("O12").Value = 0
$ws1.Range("P12").Value = 1
$ws1.Range("Q12").Value = 1
$ws1.Range("R12").Value = 0
$ws1.Range("S12").Value = 1
$ws1.Range("T12").Value = 1
$ws1.Range("U12").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 13
$ws1.Range("A13").Value = 'K:\FDTD\S2F_10_12_14.xf\Simulations\000001\Run0001\output'
$ws1.Range("B13").Value = 'S2F design +x direction 400 to 1300 nm'
$ws1.Range("C13").Value = 2.2
$ws1.Range("D13").Value = 25000
$ws1.Range("E13").Value = [double]"6.7392500000000003E-17"
$ws1.Range("F13").Value = 0
$ws1.Range("G13").Value = 400
$ws1.Range("H13").Value = 1300
$ws1.Range("I13").Value = 0
$ws1.Range("J13").Value = 2
$ws1.Range("K13").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response\Free Space 400 to 1300 5 nm mesh.mat'
$ws1.Range("L13").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response'
$ws1.Range("M13").Value = 'MultiPoint__X_sensor_4_transient_'
$ws1.Range("N13").Value = 'MultiPoint__X_sensor_4_transient_'
$ws1.Range("O13").Value = 0
$ws1.Range("P13").Value = 1
$ws1.Range("Q13").Value = 1
$ws1.Range("R13").Value = 0
$ws1.Range("S13").Value = 1
$ws1.Range("T13").Value = 1
$ws1.Range("U13").Value = 'MultiPoint__X_sensor_4_info.bin'

# New row 14
$ws1.Range("A14").Value = 'N:\Kat FDTD Data\SunPower\Cu metalization confidential full size one finger.xf\Simulations\000005\Run0001\output'
$ws1.Range("B14").Value = 'Copper back one finger full size 10-14-14 mesh 90 nm'
$ws1.Range("C14").Value = 1
$ws1.Range("D14").Value = 100000
$ws1.Range("E14").Value = [double]"1.5405899999999999E-16"
$ws1.Range("F14").Value = 0
$ws1.Range("G14").Value = 1000
$ws1.Range("H14").Value = 1300
$ws1.Range("I14").Value = 0
$ws1.Range("J14").Value = 2
$ws1.Range("K14").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response\Free Space 1000 to 1300 nm 5 nm mesh'
$ws1.Range("L14").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response'
$ws1.Range("M14").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws1.Range("N14").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws1.Range("O14").Value = 1
$ws1.Range("P14").Value = 1
$ws1.Range("Q14").Value = 0
$ws1.Range("R14").Value = 1
$ws1.Range("S14").Value = 1
$ws1.Range("T14").Value = 0
$ws1.Range("U14").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 15
$ws1.Range("A15").Value = 'M:\FDTD\Adam\post_array.xf\Simulations\000001\Run0001\output'
$ws1.Range("B15").Value = 'Adam''s post array'
$ws1.Range("C15").Value = 1.5
$ws1.Range("D15").Value = 7200
$ws1.Range("E15").Value = [double]"1.9258299999999999E-17"
$ws1.Range("F15").Value = 0
$ws1.Range("G15").Value = 400
$ws1.Range("H15").Value = 1300
$ws1.Range("I15").Value = 160
$ws1.Range("J15").Value = 160
$ws1.Range("K15").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response\Free Space 400 to 1300 5 nm mesh.mat'
$ws1.Range("L15").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response'
$ws1.Range("M15").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws1.Range("N15").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws1.Range("O15").Value = 1
$ws1.Range("P15").Value = 1
$ws1.Range("Q15").Value = 0
$ws1.Range("R15").Value = 1
$ws1.Range("S15").Value = 1
$ws1.Range("T15").Value = 0
$ws1.Range("U15").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 16
$ws1.Range("A16").Value = 'M:\FDTD\Adam\silver_tri_700_1200.xf\Simulations\000001\Run0001\output'
$ws1.Range("B16").Value = 'Mirror grating waveguide glass 1000 nm input sinusoid'
$ws1.Range("C16").Value = 2
$ws1.Range("D16").Value = 50000
$ws1.Range("E16").Value = [double]"4.8145800000000002E-17"
$ws1.Range("F16").Value = 0
$ws1.Range("G16").Value = 1000
$ws1.Range("H16").Value = 1000
$ws1.Range("I16").Value = 104
$ws1.Range("J16").Value = 4
$ws1.Range("K16").Value = 'M:\FDTD\Adam\Analyzed data\Frequency Response\sinusoidal freespace at 1 micron.mat'
$ws1.Range("L16").Value = 'M:\FDTD\Adam\Analyzed data\Frequency Response'
$ws1.Range("M16").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws1.Range("N16").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws1.Range("O16").Value = 1
$ws1.Range("P16").Value = 1
$ws1.Range("Q16").Value = 0
$ws1.Range("R16").Value = 1
$ws1.Range("S16").Value = 1
$ws1.Range("T16").Value = 0
$ws1.Range("U16").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 17
$ws1.Range("A17").Value = 'M:\FDTD\Adam\silver_tri_700_1200.xf\Simulations\000002\Run0001\output'
$ws1.Range("B17").Value = 'Mirror grating waveguide glass 1000 nm input sinusoid'
$ws1.Range("C17").Value = 2
$ws1.Range("D17").Value = 10000
$ws1.Range("E17").Value = [double]"4.8145800000000002E-17"
$ws1.Range("F17").Value = 0
$ws1.Range("G17").Value = 1000
$ws1.Range("H17").Value = 1000
$ws1.Range("I17").Value = 104
$ws1.Range("J17").Value = 4
$ws1.Range("K17").Value = 'M:\FDTD\Adam\Analyzed data\Frequency Response\sinusoidal freespace at 1 micron.mat'
$ws1.Range("L17").Value = 'M:\FDTD\Adam\Analyzed data\Frequency Response'
$ws1.Range("M17").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws1.Range("N17").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws1.Range("O17").Value = 1
$ws1.Range("P17").Value = 1
$ws1.Range("Q17").Value = 0
$ws1.Range("R17").Value = 1
$ws1.Range("S17").Value = 1
$ws1.Range("T17").Value = 0
$ws1.Range("U17").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 18
$ws1.Range("A18").Value = 'K:\FDTD\Cu confidential one finger with top texture Si-SiN.xf\Simulations\000001\Run0001\output'
$ws1.Range("B18").Value = 'Cu one finger with top texture SiN-Si'
$ws1.Range("C18").Value = 1
$ws1.Range("D18").Value = 153900
$ws1.Range("E18").Value = [double]"1.07842E-16"
$ws1.Range("F18").Value = 1
$ws1.Range("G18").Value = 1000
$ws1.Range("H18").Value = 1300
$ws1.Range("I18").Value = 4000
$ws1.Range("J18").Value = 2
$ws1.Range("K18").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response\Free Space 1000 to 1300 nm 5 nm mesh'
$ws1.Range("L18").Value = 'S:\Analyzed data\Frequency Response'
$ws1.Range("M18").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws1.Range("N18").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws1.Range("O18").Value = 1
$ws1.Range("P18").Value = 1
$ws1.Range("Q18").Value = 0
$ws1.Range("R18").Value = 1
$ws1.Range("S18").Value = 1
$ws1.Range("T18").Value = 0
$ws1.Range("U18").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 19
$ws1.Range("A19").Value = 'S:\Cu confidential one finger with top texture SiN-Si small approximation.xf\Simulations\000001\Run0001\output'
$ws1.Range("B19").Value = 'Cu one finger with top texture SiN-Si Small Approximation'
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 44600
$ws1.Range("E19").Value = [double]"1.0780700000000001E-16"
$ws1.Range("F19").Value = 0
$ws1.Range("G19").Value = 1000
$ws1.Range("H19").Value = 1300
$ws1.Range("I19").Value = 1000
$ws1.Range("J19").Value = 2
$ws1.Range("K19").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response\Free Space 1000 to 1300 nm 5 nm mesh'
$ws1.Range("L19").Value = 'S:\Analyzed data\Frequency Response'
$ws1.Range("M19").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws1.Range("N19").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws1.Range("O19").Value = 1
$ws1.Range("P19").Value = 1
$ws1.Range("Q19").Value = 0
$ws1.Range("R19").Value = 1
$ws1.Range("S19").Value = 1
$ws1.Range("T19").Value = 0
$ws1.Range("U19").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# Column U width (best-fit-ish) and freeze-pane selection restore
$ws1.Range("U1").ColumnWidth = 37
$ws1.Activate()
$ws1.Range("F20").Select()

# ---------------------------------------------------------------------------
# Sheet "Free Space Raw FDTD Data"
# ---------------------------------------------------------------------------

# New header columns N:T (field coefficients + info file stub)
$ws2.Range("N1").Value = 'Ex Coefficient'
$ws2.Range("O1").Value = 'Ey Coefficient'
$ws2.Range("P1").Value = 'Ez Coefficient'
$ws2.Range("Q1").Value = 'Hx Coefficient'
$ws2.Range("R1").Value = 'Hy Coefficient'
$ws2.Range("S1").Value = 'Hz Coefficient'
$ws2.Range("T1").Value = 'Info file stub'
$ws2.Range("A1").EntireRow.RowHeight = 43.2

# Existing rows 2-3 get the new coefficient + info values
$ws2.Range("N2").Value = 1
$ws2.Range("O2").Value = 1
$ws2.Range("P2").Value = 0
$ws2.Range("Q2").Value = 1
$ws2.Range("R2").Value = 1
$ws2.Range("S2").Value = 0
$ws2.Range("T2").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

$ws2.Range("N3").Value = 1
$ws2.Range("O3").Value = 1
$ws2.Range("P3").Value = 0
$ws2.Range("Q3").Value = 1
$ws2.Range("R3").Value = 1
$ws2.Range("S3").Value = 0
$ws2.Range("T3").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 4
$ws2.Range("A4").Value = 'K:\FDTD\Free space 400 to 1300.xf\Simulations\000001\Run0001\output'
$ws2.Range("B4").Value = 'Free Space 400 to 1300 5 nm mesh'
$ws2.Range("C4").Value = 1
$ws2.Range("D4").Value = 1700
$ws2.Range("E4").Value = [double]"9.6291700000000004E-18"
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 400
$ws2.Range("H4").Value = 1300
$ws2.Range("I4").Value = 200
$ws2.Range("J4").Value = 2
$ws2.Range("K4").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response'
$ws2.Range("L4").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M4").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N4").Value = 1
$ws2.Range("O4").Value = 1
$ws2.Range("P4").Value = 0
$ws2.Range("Q4").Value = 1
$ws2.Range("R4").Value = 1
$ws2.Range("S4").Value = 0
$ws2.Range("T4").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 5
$ws2.Range("A5").Value = 'N:\Kat FDTD Data\SunPower\Free Space 1000 to 1300 nm.xf\Simulations\000001\Run0001\output'
$ws2.Range("B5").Value = 'Free Space 1000 to 1300 nm 5 nm mesh'
$ws2.Range("C5").Value = 1
$ws2.Range("D5").Value = 9600
$ws2.Range("E5").Value = [double]"9.6291700000000004E-18"
$ws2.Range("F5").Value = 0
$ws2.Range("G5").Value = 1000
$ws2.Range("H5").Value = 1300
$ws2.Range("I5").Value = 200
$ws2.Range("J5").Value = 2
$ws2.Range("K5").Value = 'N:\Kat FDTD Data\Analyzed Data\Frequency Response'
$ws2.Range("L5").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M5").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N5").Value = 1
$ws2.Range("O5").Value = 1
$ws2.Range("P5").Value = 0
$ws2.Range("Q5").Value = 1
$ws2.Range("R5").Value = 1
$ws2.Range("S5").Value = 0
$ws2.Range("T5").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 6
$ws2.Range("A6").Value = 'M:\FDTD\Adam\freespace_sinusoid_1micron.xf\Simulations\000001\Run0001\output'
$ws2.Range("B6").Value = 'sinusoidal freespace at 1 micron'
$ws2.Range("C6").Value = 1
$ws2.Range("D6").Value = 3000
$ws2.Range("E6").Value = [double]"4.8145800000000002E-17"
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 1000
$ws2.Range("H6").Value = 1000
$ws2.Range("I6").Value = 104
$ws2.Range("J6").Value = 4
$ws2.Range("K6").Value = 'M:\FDTD\Adam\Analyzed data\Frequency Response'
$ws2.Range("L6").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M6").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N6").Value = 1
$ws2.Range("O6").Value = 1
$ws2.Range("P6").Value = 0
$ws2.Range("Q6").Value = 1
$ws2.Range("R6").Value = 1
$ws2.Range("S6").Value = 0
$ws2.Range("T6").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 7
$ws2.Range("A7").Value = 'S:\Free Space.xf\Simulations\000002\Run0001\output'
$ws2.Range("B7").Value = 'Free Space 400-500 nm'
$ws2.Range("C7").Value = 1
$ws2.Range("D7").Value = 4700
$ws2.Range("E7").Value = [double]"9.6291700000000004E-18"
$ws2.Range("F7").Value = 0
$ws2.Range("G7").Value = 400
$ws2.Range("H7").Value = 500
$ws2.Range("I7").Value = 200
$ws2.Range("J7").Value = 2
$ws2.Range("K7").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L7").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M7").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N7").Value = 1
$ws2.Range("O7").Value = 1
$ws2.Range("P7").Value = 0
$ws2.Range("Q7").Value = 1
$ws2.Range("R7").Value = 1
$ws2.Range("S7").Value = 0
$ws2.Range("T7").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 8
$ws2.Range("A8").Value = 'S:\Free Space.xf\Simulations\000002\Run0002\output'
$ws2.Range("B8").Value = 'Free Space 500-600 nm'
$ws2.Range("C8").Value = 1
$ws2.Range("D8").Value = 6800
$ws2.Range("E8").Value = [double]"9.6291700000000004E-18"
$ws2.Range("F8").Value = 0
$ws2.Range("G8").Value = 500
$ws2.Range("H8").Value = 600
$ws2.Range("I8").Value = 200
$ws2.Range("J8").Value = 2
$ws2.Range("K8").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L8").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M8").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N8").Value = 1
$ws2.Range("O8").Value = 1
$ws2.Range("P8").Value = 0
$ws2.Range("Q8").Value = 1
$ws2.Range("R8").Value = 1
$ws2.Range("S8").Value = 0
$ws2.Range("T8").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 9
$ws2.Range("A9").Value = 'S:\Free Space.xf\Simulations\000002\Run0003\output'
$ws2.Range("B9").Value = 'Free Space 600-700 nm'
$ws2.Range("C9").Value = 1
$ws2.Range("D9").Value = 9300
$ws2.Range("E9").Value = [double]"9.6291700000000004E-18"
$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = 600
$ws2.Range("H9").Value = 700
$ws2.Range("I9").Value = 200
$ws2.Range("J9").Value = 2
$ws2.Range("K9").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L9").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M9").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N9").Value = 1
$ws2.Range("O9").Value = 1
$ws2.Range("P9").Value = 0
$ws2.Range("Q9").Value = 1
$ws2.Range("R9").Value = 1
$ws2.Range("S9").Value = 0
$ws2.Range("T9").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 10
$ws2.Range("A10").Value = 'S:\Free Space.xf\Simulations\000004\Run0001\output'
$ws2.Range("B10").Value = 'Free Space 700-800 nm'
$ws2.Range("C10").Value = 1
$ws2.Range("D10").Value = 6200
$ws2.Range("E10").Value = [double]"1.9258299999999999E-17"
$ws2.Range("F10").Value = 0
$ws2.Range("G10").Value = 700
$ws2.Range("H10").Value = 800
$ws2.Range("I10").Value = 100
$ws2.Range("J10").Value = 2
$ws2.Range("K10").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L10").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M10").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N10").Value = 1
$ws2.Range("O10").Value = 1
$ws2.Range("P10").Value = 0
$ws2.Range("Q10").Value = 1
$ws2.Range("R10").Value = 1
$ws2.Range("S10").Value = 0
$ws2.Range("T10").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 11
$ws2.Range("A11").Value = 'S:\Free Space.xf\Simulations\000004\Run0002\output'
$ws2.Range("B11").Value = 'Free Space 800-900 nm'
$ws2.Range("C11").Value = 1
$ws2.Range("D11").Value = 7900
$ws2.Range("E11").Value = [double]"1.9258299999999999E-17"
$ws2.Range("F11").Value = 1
$ws2.Range("G11").Value = 800
$ws2.Range("H11").Value = 900
$ws2.Range("I11").Value = 100
$ws2.Range("J11").Value = 2
$ws2.Range("K11").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L11").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M11").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N11").Value = 1
$ws2.Range("O11").Value = 1
$ws2.Range("P11").Value = 0
$ws2.Range("Q11").Value = 1
$ws2.Range("R11").Value = 1
$ws2.Range("S11").Value = 0
$ws2.Range("T11").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 12
$ws2.Range("A12").Value = 'S:\Free Space.xf\Simulations\000004\Run0003\output'
$ws2.Range("B12").Value = 'Free Space 900-1000 nm'
$ws2.Range("C12").Value = 1
$ws2.Range("D12").Value = 9900
$ws2.Range("E12").Value = [double]"1.9258299999999999E-17"
$ws2.Range("F12").Value = 1
$ws2.Range("G12").Value = 900
$ws2.Range("H12").Value = 1000
$ws2.Range("I12").Value = 100
$ws2.Range("J12").Value = 2
$ws2.Range("K12").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L12").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M12").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N12").Value = 1
$ws2.Range("O12").Value = 1
$ws2.Range("P12").Value = 0
$ws2.Range("Q12").Value = 1
$ws2.Range("R12").Value = 1
$ws2.Range("S12").Value = 0
$ws2.Range("T12").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 13
$ws2.Range("A13").Value = 'S:\Free Space.xf\Simulations\000004\Run0004\output'
$ws2.Range("B13").Value = 'Free Space 1000-1100 nm'
$ws2.Range("C13").Value = 1
$ws2.Range("D13").Value = 12000
$ws2.Range("E13").Value = [double]"1.9258299999999999E-17"
$ws2.Range("F13").Value = 1
$ws2.Range("G13").Value = 1000
$ws2.Range("H13").Value = 1100
$ws2.Range("I13").Value = 100
$ws2.Range("J13").Value = 2
$ws2.Range("K13").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L13").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M13").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N13").Value = 1
$ws2.Range("O13").Value = 1
$ws2.Range("P13").Value = 0
$ws2.Range("Q13").Value = 1
$ws2.Range("R13").Value = 1
$ws2.Range("S13").Value = 0
$ws2.Range("T13").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 14
$ws2.Range("A14").Value = 'S:\Free Space.xf\Simulations\000004\Run0005\output'
$ws2.Range("B14").Value = 'Free Space 1100-1200 nm'
$ws2.Range("C14").Value = 1
$ws2.Range("D14").Value = 14300
$ws2.Range("E14").Value = [double]"1.9258299999999999E-17"
$ws2.Range("F14").Value = 1
$ws2.Range("G14").Value = 1100
$ws2.Range("H14").Value = 1200
$ws2.Range("I14").Value = 100
$ws2.Range("J14").Value = 2
$ws2.Range("K14").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L14").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M14").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N14").Value = 1
$ws2.Range("O14").Value = 1
$ws2.Range("P14").Value = 0
$ws2.Range("Q14").Value = 1
$ws2.Range("R14").Value = 1
$ws2.Range("S14").Value = 0
$ws2.Range("T14").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

# New row 15
$ws2.Range("A15").Value = 'S:\Free Space.xf\Simulations\000004\Run0006\output'
$ws2.Range("B15").Value = 'Free Space 1200-1300 nm'
$ws2.Range("C15").Value = 1
$ws2.Range("D15").Value = 16900
$ws2.Range("E15").Value = [double]"1.9258299999999999E-17"
$ws2.Range("F15").Value = 1
$ws2.Range("G15").Value = 1200
$ws2.Range("H15").Value = 1300
$ws2.Range("I15").Value = 100
$ws2.Range("J15").Value = 2
$ws2.Range("K15").Value = 'S:\Analyzed Data\Frequency Response'
$ws2.Range("L15").Value = 'MultiPoint_Reflected_Light_Top_0_transient_'
$ws2.Range("M15").Value = 'MultiPoint_Transmitted_Light_Bottom_1_transient_'
$ws2.Range("N15").Value = 1
$ws2.Range("O15").Value = 1
$ws2.Range("P15").Value = 0
$ws2.Range("Q15").Value = 1
$ws2.Range("R15").Value = 1
$ws2.Range("S15").Value = 0
$ws2.Range("T15").Value = 'MultiPoint_Reflected_Light_Top_0_info.bin'

$ws2.Activate()
$ws2.Range("G5").Select()
